$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "ZSO-0025"

$ws.Range("B3").Value = "DSR-0350"
$ws.Range("C3").Value = "Mayer Doa Varaitey Store"
$ws.Range("D3").Value = "Lalpur"
$ws.Range("E3").Value = "Md Abu Salek"
$ws.Range("I3").Value = "Md Abu Salek"
$ws.Range("K3").Value = "Natore"
$ws.Range("L3").Value = "Lalpur"
$ws.Range("M3").Value = "ZSO-0025"
$ws.Range("N3").Value = "Chinir Bottola, Lalpur, Natore"
$ws.Range("P3").Value = 1738027070
$ws.Range("T3").Value = 1738027070

$ws.Range("I19").Select()
